# Fixed most sound issues
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fix typo "Double  Vowel" (double space) -> "Double Vowel" (single space)
# in the Female Audio column (E) for the double-vowel rows.
$doubleVowelRows = @(51, 52, 53, 54, 55, 56)
foreach ($r in $doubleVowelRows) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value = $cell.Value2.Replace("Double  Vowel", "Double Vowel")
}

# --- Fix "Female - Consonant X.wav" -> "Female - Single Consonant X.wav"
# for the consonant rows that were missing "Single" in the Female Audio column (E).
$consonantFixRows = @(58, 60, 62, 64, 66, 69)
foreach ($r in $consonantFixRows) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value = $cell.Value2.Replace("Female - Consonant", "Female - Single Consonant")
}

# --- Set the Pronunciation Video column (F) to "y" for the other consonant rows.
$fMarkRows = @(57, 59, 61, 63, 65, 67, 68, 70, 71, 72, 73)
foreach ($r in $fMarkRows) {
    $ws.Cells.Item($r, 6).Value = "y"
}

# --- Touch C74 so it carries the same formatting as the rows above it.
$ws.Range("C73").Copy()
$ws.Range("C74").PasteSpecial(-4122)
$ws.Range("C74").ClearContents()
$excel.CutCopyMode = 0

# --- Restore the selection / scroll position left behind by the edits.
$ws.Range("E61").Select()
$excel.ActiveWindow.ScrollRow = 62
$excel.ActiveWindow.ScrollColumn = 2
